# Auto-generated edit script: apply value updates to Sheets per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 1112167.2
$ws.Cells.Item(9, 9).Value = 1667816.9
$ws.Cells.Item(9, 10).Value = 868
$ws.Cells.Item(9, 11).Value = 1667816.9
$ws.Cells.Item(9, 12).Value = 868
$ws.Cells.Item(9, 13).Value = -1667647.9
$ws.Cells.Item(9, 14).Value = -1206
$ws.Cells.Item(15, 8).Value = 573.881
$ws.Cells.Item(15, 9).Value = 573.881
$ws.Cells.Item(15, 11).Value = 1721.643
$ws.Cells.Item(15, 13).Value = -1552.643
$ws.Cells.Item(32, 8).Value = 5748.375
$ws.Cells.Item(32, 9).Value = 7650
$ws.Cells.Item(32, 10).Value = 5476.7144
$ws.Cells.Item(32, 11).Value = 7650
$ws.Cells.Item(32, 12).Value = 5476.7144
$ws.Cells.Item(32, 13).Value = -7324
$ws.Cells.Item(32, 14).Value = -6128.7144
$ws.Cells.Item(38, 8).Value = 8760.549999999999
$ws.Cells.Item(38, 9).Value = 11721.8
$ws.Cells.Item(38, 11).Value = 35165.39999999999
$ws.Cells.Item(38, 13).Value = -34793.39999999999
$ws.Cells.Item(40, 8).Value = 3510.3333
$ws.Cells.Item(40, 9).Value = 3250
$ws.Cells.Item(40, 10).Value = 3562.4
$ws.Cells.Item(40, 11).Value = 3250
$ws.Cells.Item(40, 12).Value = 3562.4
$ws.Cells.Item(40, 13).Value = -3075
$ws.Cells.Item(40, 14).Value = -3912.4
$ws.Cells.Item(43, 8).Value = 7333.278
$ws.Cells.Item(43, 10).Value = 7999.9375
$ws.Cells.Item(43, 12).Value = 7999.9375
$ws.Cells.Item(43, 14).Value = -8137.9375
$ws.Cells.Item(46, 8).Value = 55833.168
$ws.Cells.Item(46, 10).Value = 49999.5
$ws.Cells.Item(46, 12).Value = 149998.5
$ws.Cells.Item(46, 14).Value = -150236.5
$ws.Cells.Item(53, 8).Value = 371.75
$ws.Cells.Item(53, 9).Value = 220.75
$ws.Cells.Item(53, 11).Value = 220.75
$ws.Cells.Item(53, 13).Value = 416.25
$ws.Cells.Item(58, 8).Value = 50000070
$ws.Cells.Item(58, 9).Value = 50000070
$ws.Cells.Item(58, 11).Value = 150000210
$ws.Cells.Item(58, 13).Value = -150000060
$ws.Cells.Item(60, 8).Value = 55833.168
$ws.Cells.Item(60, 10).Value = 49999.5
$ws.Cells.Item(60, 12).Value = 149998.5
$ws.Cells.Item(60, 14).Value = -150966.5
$ws.Cells.Item(74, 8).Value = 7522.1816
$ws.Cells.Item(74, 9).Value = 5304.8887
$ws.Cells.Item(74, 11).Value = 5304.8887
$ws.Cells.Item(74, 13).Value = -4368.8887
$ws.Cells.Item(77, 8).Value = 7522.1816
$ws.Cells.Item(77, 9).Value = 5304.8887
$ws.Cells.Item(77, 11).Value = 26524.4435
$ws.Cells.Item(77, 13).Value = -21844.4435
$ws.Cells.Item(95, 8).Value = 0
$ws.Cells.Item(95, 10).Value = 0
$ws.Cells.Item(95, 12).Value = 0
$ws.Cells.Item(95, 14).ClearContents()
$ws.Cells.Item(96, 8).Value = 1111503
$ws.Cells.Item(96, 10).Value = 397.5
$ws.Cells.Item(96, 12).Value = 1192.5
$ws.Cells.Item(96, 14).Value = -3938.5
$ws.Cells.Item(104, 8).Value = 201.28572
$ws.Cells.Item(104, 9).Value = 201.28572
$ws.Cells.Item(104, 11).Value = 603.85716
$ws.Cells.Item(104, 13).Value = 1143.14284
$ws.Cells.Item(115, 8).Value = 184
$ws.Cells.Item(115, 9).Value = 184
$ws.Cells.Item(115, 11).Value = 552
$ws.Cells.Item(115, 13).Value = 1015
$ws.Cells.Item(129, 8).Value = 1269.5
$ws.Cells.Item(129, 10).Value = 2153
$ws.Cells.Item(129, 12).Value = 6459
$ws.Cells.Item(129, 14).Value = -16459
$ws.Cells.Item(137, 8).Value = 11829.0625
$ws.Cells.Item(137, 9).Value = 15396.909
$ws.Cells.Item(137, 10).Value = 3979.8
$ws.Cells.Item(137, 11).Value = 46190.727
$ws.Cells.Item(137, 12).Value = 11939.4
$ws.Cells.Item(137, 13).Value = -43640.727
$ws.Cells.Item(137, 14).Value = -17039.4
$ws.Cells.Item(138, 8).Value = 3324.9805
$ws.Cells.Item(138, 9).Value = 2320.6
$ws.Cells.Item(138, 10).Value = 3743.4722
$ws.Cells.Item(138, 11).Value = 6961.799999999999
$ws.Cells.Item(138, 12).Value = 11230.4166
$ws.Cells.Item(138, 13).Value = -1821.799999999999
$ws.Cells.Item(138, 14).Value = -21510.4166

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 4287.6665
$ws.Cells.Item(45, 9).Value = 1518.8
$ws.Cells.Item(45, 10).Value = 7748.75
$ws.Cells.Item(45, 11).Value = 1518.8
$ws.Cells.Item(45, 12).Value = 7748.75
$ws.Cells.Item(45, 13).Value = -1141.8
$ws.Cells.Item(45, 14).Value = -8502.75
$ws.Cells.Item(61, 8).Value = 4624.04
$ws.Cells.Item(61, 9).Value = 5153.4116
$ws.Cells.Item(61, 10).Value = 3499.125
$ws.Cells.Item(61, 11).Value = 5153.4116
$ws.Cells.Item(61, 12).Value = 3499.125
$ws.Cells.Item(61, 13).Value = -4941.4116
$ws.Cells.Item(61, 14).Value = -3923.125
$ws.Cells.Item(74, 8).Value = 3504.9583
$ws.Cells.Item(74, 9).Value = 3748.5293
$ws.Cells.Item(74, 10).Value = 2913.4285
$ws.Cells.Item(74, 11).Value = 3748.5293
$ws.Cells.Item(74, 12).Value = 2913.4285
$ws.Cells.Item(74, 13).Value = -2874.5293
$ws.Cells.Item(74, 14).Value = -4661.4285
$ws.Cells.Item(77, 8).Value = 3504.9583
$ws.Cells.Item(77, 9).Value = 3748.5293
$ws.Cells.Item(77, 10).Value = 2913.4285
$ws.Cells.Item(77, 11).Value = 18742.6465
$ws.Cells.Item(77, 12).Value = 14567.1425
$ws.Cells.Item(77, 13).Value = -14374.6465
$ws.Cells.Item(77, 14).Value = -23303.1425
$ws.Cells.Item(132, 8).Value = 4075
$ws.Cells.Item(132, 9).Value = 4075
$ws.Cells.Item(132, 11).Value = 12225
$ws.Cells.Item(132, 13).Value = -9695
$ws.Cells.Item(136, 8).Value = 4624.04
$ws.Cells.Item(136, 9).Value = 5153.4116
$ws.Cells.Item(136, 10).Value = 3499.125
$ws.Cells.Item(136, 11).Value = 15460.2348
$ws.Cells.Item(136, 12).Value = 10497.375
$ws.Cells.Item(136, 13).Value = -12910.2348
$ws.Cells.Item(136, 14).Value = -15597.375
$ws.Cells.Item(138, 8).Value = 120000
$ws.Cells.Item(138, 10).Value = 120000
$ws.Cells.Item(138, 12).Value = 120000
$ws.Cells.Item(138, 14).Value = -130280

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 97110.836
$ws.Cells.Item(94, 9).Value = 120348.17
$ws.Cells.Item(94, 11).Value = 120348.17
$ws.Cells.Item(94, 13).Value = -119897.17
$ws.Cells.Item(101, 8).Value = 0
$ws.Cells.Item(101, 10).Value = 0
$ws.Cells.Item(101, 12).Value = 0
$ws.Cells.Item(101, 14).ClearContents()
$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 14).ClearContents()
$ws.Cells.Item(103, 8).Value = 14833.167
$ws.Cells.Item(103, 10).Value = 14833.167
$ws.Cells.Item(103, 12).Value = 14833.167
$ws.Cells.Item(103, 14).Value = -17177.167
$ws.Cells.Item(134, 8).Value = 3178.5833
$ws.Cells.Item(134, 9).Value = 3204.8333
$ws.Cells.Item(134, 11).Value = 9614.499899999999
$ws.Cells.Item(134, 13).Value = -7079.499899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3031.2222
$ws.Cells.Item(31, 9).Value = 1752.9166
$ws.Cells.Item(31, 10).Value = 4053.8667
$ws.Cells.Item(31, 11).Value = 1752.9166
$ws.Cells.Item(31, 12).Value = 4053.8667
$ws.Cells.Item(31, 13).Value = -1457.9166
$ws.Cells.Item(31, 14).Value = -4643.8667
$ws.Cells.Item(34, 8).Value = 3031.2222
$ws.Cells.Item(34, 9).Value = 1752.9166
$ws.Cells.Item(34, 10).Value = 4053.8667
$ws.Cells.Item(34, 11).Value = 1752.9166
$ws.Cells.Item(34, 12).Value = 4053.8667
$ws.Cells.Item(34, 13).Value = -1550.9166
$ws.Cells.Item(34, 14).Value = -4457.8667
$ws.Cells.Item(58, 8).Value = 1910.9333
$ws.Cells.Item(58, 9).Value = 1763.6666
$ws.Cells.Item(58, 10).Value = 2500
$ws.Cells.Item(58, 11).Value = 1763.6666
$ws.Cells.Item(58, 12).Value = 2500
$ws.Cells.Item(58, 13).Value = -1560.6666
$ws.Cells.Item(58, 14).Value = -2906
$ws.Cells.Item(102, 8).Value = 35000
$ws.Cells.Item(102, 9).Value = 35000
$ws.Cells.Item(102, 11).Value = 35000
$ws.Cells.Item(102, 13).Value = -32566
$ws.Cells.Item(104, 8).Value = 49985
$ws.Cells.Item(104, 10).Value = 49985
$ws.Cells.Item(104, 12).Value = 49985
$ws.Cells.Item(104, 14).Value = -55227
$ws.Cells.Item(132, 8).Value = 2152
$ws.Cells.Item(132, 9).Value = 2196
$ws.Cells.Item(132, 11).Value = 6588
$ws.Cells.Item(132, 13).Value = -4058
$ws.Cells.Item(134, 8).Value = 1531.6471
$ws.Cells.Item(134, 9).Value = 1467.3334
$ws.Cells.Item(134, 11).Value = 4402.0002
$ws.Cells.Item(134, 13).Value = -1867.0002
$ws.Cells.Item(135, 8).Value = 71827.25
$ws.Cells.Item(135, 10).Value = 71827.25
$ws.Cells.Item(135, 12).Value = 71827.25
$ws.Cells.Item(135, 14).Value = -81967.25
$ws.Cells.Item(136, 8).Value = 1910.9333
$ws.Cells.Item(136, 9).Value = 1763.6666
$ws.Cells.Item(136, 10).Value = 2500
$ws.Cells.Item(136, 11).Value = 5290.9998
$ws.Cells.Item(136, 12).Value = 7500
$ws.Cells.Item(136, 13).Value = -2740.9998
$ws.Cells.Item(136, 14).Value = -12600
$ws.Cells.Item(141, 8).Value = 51499.668
$ws.Cells.Item(141, 10).Value = 94499
$ws.Cells.Item(141, 12).Value = 94499
$ws.Cells.Item(141, 14).Value = -104859

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 27181318
$ws.Cells.Item(4, 10).Value = 7562648.5
$ws.Cells.Item(4, 12).Value = 22687945.5
$ws.Cells.Item(4, 14).Value = -22688169.5
$ws.Cells.Item(22, 8).Value = 415
$ws.Cells.Item(22, 10).Value = 415
$ws.Cells.Item(22, 12).Value = 1245
$ws.Cells.Item(22, 14).Value = -1583
$ws.Cells.Item(27, 8).Value = 415
$ws.Cells.Item(27, 10).Value = 415
$ws.Cells.Item(27, 12).Value = 1245
$ws.Cells.Item(27, 14).Value = -1449
$ws.Cells.Item(81, 8).Value = 8049.75
$ws.Cells.Item(81, 10).Value = 8327.182000000001
$ws.Cells.Item(81, 12).Value = 24981.546
$ws.Cells.Item(81, 14).Value = -27227.546
$ws.Cells.Item(84, 8).Value = 8049.75
$ws.Cells.Item(84, 10).Value = 8327.182000000001
$ws.Cells.Item(84, 12).Value = 74944.63800000001
$ws.Cells.Item(84, 14).Value = -86176.63800000001
$ws.Cells.Item(96, 8).Value = 10000
$ws.Cells.Item(96, 10).Value = 10000
$ws.Cells.Item(96, 12).Value = 30000
$ws.Cells.Item(96, 14).Value = -34118
$ws.Cells.Item(132, 8).Value = 1981.7059
$ws.Cells.Item(132, 9).Value = 1333.3334
$ws.Cells.Item(132, 10).Value = 2120.6428
$ws.Cells.Item(132, 11).Value = 12000.0006
$ws.Cells.Item(132, 12).Value = 19085.7852
$ws.Cells.Item(132, 13).Value = -9470.000599999999
$ws.Cells.Item(132, 14).Value = -24145.7852
$ws.Cells.Item(134, 8).Value = 2169.7144
$ws.Cells.Item(137, 8).Value = 3376.4
$ws.Cells.Item(137, 9).Value = 3279.625
$ws.Cells.Item(137, 10).Value = 3763.5
$ws.Cells.Item(137, 11).Value = 9838.875
$ws.Cells.Item(137, 12).Value = 11290.5
$ws.Cells.Item(137, 13).Value = -4738.875
$ws.Cells.Item(137, 14).Value = -21490.5
$ws.Cells.Item(139, 8).Value = 4003.75
$ws.Cells.Item(139, 9).Value = 2209
$ws.Cells.Item(139, 10).Value = 5285.7144
$ws.Cells.Item(139, 11).Value = 6627
$ws.Cells.Item(139, 12).Value = 15857.1432
$ws.Cells.Item(139, 13).Value = -1487
$ws.Cells.Item(139, 14).Value = -26137.1432
$ws.Cells.Item(141, 8).Value = 2165
$ws.Cells.Item(141, 9).Value = 2165
$ws.Cells.Item(141, 10).Value = 0
$ws.Cells.Item(141, 11).Value = 6495
$ws.Cells.Item(141, 12).Value = 0
$ws.Cells.Item(141, 13).Value = -1315
$ws.Cells.Item(141, 14).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(12, 8).Value = 531787.4399999999
$ws.Cells.Item(12, 10).Value = 1264999.6
$ws.Cells.Item(12, 12).Value = 1264999.6
$ws.Cells.Item(12, 14).Value = -1265279.6
$ws.Cells.Item(26, 8).Value = 34875
$ws.Cells.Item(26, 10).Value = 34875
$ws.Cells.Item(26, 12).Value = 34875
$ws.Cells.Item(26, 14).Value = -35435
$ws.Cells.Item(50, 8).Value = 34875
$ws.Cells.Item(50, 10).Value = 34875
$ws.Cells.Item(50, 12).Value = 34875
$ws.Cells.Item(50, 14).Value = -35871
$ws.Cells.Item(70, 8).Value = 6544736.5
$ws.Cells.Item(70, 9).Value = 22230428
$ws.Cells.Item(70, 11).Value = 22230428
$ws.Cells.Item(70, 13).Value = -22230158
$ws.Cells.Item(73, 8).Value = 6544736.5
$ws.Cells.Item(73, 9).Value = 22230428
$ws.Cells.Item(73, 11).Value = 22230428
$ws.Cells.Item(73, 13).Value = -22229492
$ws.Cells.Item(80, 8).Value = 15499.454
$ws.Cells.Item(80, 9).Value = 19213.428
$ws.Cells.Item(80, 11).Value = 19213.428
$ws.Cells.Item(80, 13).Value = -18215.428
$ws.Cells.Item(83, 8).Value = 15499.454
$ws.Cells.Item(83, 9).Value = 19213.428
$ws.Cells.Item(83, 11).Value = 96067.14
$ws.Cells.Item(83, 13).Value = -91075.14
$ws.Cells.Item(97, 8).Value = 4302.2905
$ws.Cells.Item(97, 9).Value = 1194.6538
$ws.Cells.Item(97, 10).Value = 20462
$ws.Cells.Item(97, 11).Value = 1194.6538
$ws.Cells.Item(97, 12).Value = 20462
$ws.Cells.Item(97, 13).Value = -698.6538
$ws.Cells.Item(97, 14).Value = -21454
$ws.Cells.Item(132, 8).Value = 3222.5
$ws.Cells.Item(132, 9).Value = 3148.3684
$ws.Cells.Item(132, 11).Value = 9445.1052
$ws.Cells.Item(132, 13).Value = -6915.1052
$ws.Cells.Item(134, 8).Value = 71585.2
$ws.Cells.Item(134, 10).Value = 71585.2
$ws.Cells.Item(134, 12).Value = 214755.6
$ws.Cells.Item(134, 14).Value = -219825.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5067.143
$ws.Cells.Item(7, 9).Value = 5494.6
$ws.Cells.Item(7, 10).Value = 3998.5
$ws.Cells.Item(7, 11).Value = 5494.6
$ws.Cells.Item(7, 12).Value = 3998.5
$ws.Cells.Item(7, 13).Value = -5382.6
$ws.Cells.Item(7, 14).Value = -4222.5
$ws.Cells.Item(17, 8).Value = 5214.2856
$ws.Cells.Item(17, 9).Value = 2750
$ws.Cells.Item(17, 11).Value = 2750
$ws.Cells.Item(17, 13).Value = -2580
$ws.Cells.Item(22, 8).Value = 3587.6875
$ws.Cells.Item(22, 9).Value = 3207.625
$ws.Cells.Item(22, 11).Value = 3207.625
$ws.Cells.Item(22, 13).Value = -2912.625
$ws.Cells.Item(27, 8).Value = 3587.6875
$ws.Cells.Item(27, 9).Value = 3207.625
$ws.Cells.Item(27, 11).Value = 3207.625
$ws.Cells.Item(27, 13).Value = -3100.625
$ws.Cells.Item(81, 8).Value = 124995
$ws.Cells.Item(81, 10).Value = 124995
$ws.Cells.Item(81, 12).Value = 124995
$ws.Cells.Item(81, 14).Value = -126991
$ws.Cells.Item(84, 8).Value = 124995
$ws.Cells.Item(84, 10).Value = 124995
$ws.Cells.Item(84, 12).Value = 374985
$ws.Cells.Item(84, 14).Value = -384969
$ws.Cells.Item(121, 8).Value = 242000
$ws.Cells.Item(121, 10).Value = 242000
$ws.Cells.Item(121, 12).Value = 242000
$ws.Cells.Item(121, 14).Value = -245494
$ws.Cells.Item(126, 8).Value = 5067.143
$ws.Cells.Item(126, 9).Value = 5494.6
$ws.Cells.Item(126, 10).Value = 3998.5
$ws.Cells.Item(126, 11).Value = 16483.8
$ws.Cells.Item(126, 12).Value = 11995.5
$ws.Cells.Item(126, 13).Value = -14013.8
$ws.Cells.Item(126, 14).Value = -16935.5
$ws.Cells.Item(132, 8).Value = 5346.1177
$ws.Cells.Item(132, 9).Value = 5228.143
$ws.Cells.Item(132, 11).Value = 15684.429
$ws.Cells.Item(132, 13).Value = -13154.429
$ws.Cells.Item(136, 8).Value = 3829.375
$ws.Cells.Item(136, 9).Value = 3878.0667
$ws.Cells.Item(136, 11).Value = 11634.2001
$ws.Cells.Item(136, 13).Value = -9084.2001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 14).ClearContents()
$ws.Cells.Item(107, 8).Value = 1082.5625
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 14).ClearContents()
$ws.Cells.Item(111, 8).Value = 0
$ws.Cells.Item(111, 10).Value = 0
$ws.Cells.Item(111, 12).Value = 0
$ws.Cells.Item(111, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 6610.9443
$ws.Cells.Item(122, 10).Value = 7200.3335
$ws.Cells.Item(122, 12).Value = 21601.0005
$ws.Cells.Item(122, 14).Value = -26501.0005
$ws.Cells.Item(132, 8).Value = 2601.0667
$ws.Cells.Item(132, 9).Value = 2265.4285
$ws.Cells.Item(132, 11).Value = 6796.2855
$ws.Cells.Item(132, 13).Value = -4266.2855
$ws.Cells.Item(136, 8).Value = 1585.5217
$ws.Cells.Item(136, 9).Value = 1387.2106
$ws.Cells.Item(136, 10).Value = 2527.5
$ws.Cells.Item(136, 11).Value = 4161.6318
$ws.Cells.Item(136, 12).Value = 7582.5
$ws.Cells.Item(136, 13).Value = -1611.6318
$ws.Cells.Item(136, 14).Value = -12682.5
$ws.Cells.Item(140, 8).Value = 79998
$ws.Cells.Item(140, 10).Value = 79998
$ws.Cells.Item(140, 12).Value = 79998
$ws.Cells.Item(140, 14).Value = -90358

